# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Numeric-looking Price (D) and Volume (E) columns are forced to remain plain text
# (matching the source inlineStr cells) by entering them with a leading apostrophe,
# then resetting the cell style back to "Normal" so no stray number-format / quote-prefix
# style sticks around on the cell (keeps the XML diff to just the <is><t> value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'70.452.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.51%  "
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = "'3.632.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.23%  "
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = "'594.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.09%  "
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'195.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.62%  "
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("E7").Value = "'  +1.92%  "
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = "'3.625.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.21%  "
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("E9").Value = "'  -0.10%  "
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("E10").Value = "'  +2.60%  "
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("E11").Value = "'  +2.97%  "
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "'58.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.25%  "
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").Value = "'0.0000291"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.94%  "
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("E14").Value = "'  +5.15%  "
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = "'4.212.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.16%  "
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = "'19.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.03%  "
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").Value = "'3.634.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.93%  "
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = "'70.445.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.47%  "
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = "'12.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.54%  "
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("E20").Value = "'  +1.76%  "
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("E21").Value = "'  +4.44%  "
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = "'488.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.29%  "
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = "'19.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +16.13%  "
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = "'5.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.24%  "
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = "'4.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.19%  "
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").Value = "'91.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.30%  "
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'3.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +6.93%  "
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("D28").Value = "'11.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.68%  "
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("D29").Value = "'9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.94%  "
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +10.94%  "
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'33.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.61%  "
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.25%  "
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'627.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.44%  "
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("D34").Value = "'12.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.69%  "
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("D35").Value = "'66.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.41%  "
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("D36").Value = "'40.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +10.94%  "
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("D37").Value = "'0.416"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.59%  "
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("D38").Value = "'0.0₃0825"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +7.39%  "
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.86%  "
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.01%  "
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.54%  "
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = "'3.293.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.41%  "
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("E43").Value = "'  +8.20%  "
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'2.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +12.34%  "
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = "'0.0454"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.65%  "
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").Value = "'2.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.29%  "
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").Value = "'  +2.03%  "
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("E48").Value = "'  +2.56%  "
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").Value = "'9.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.50%  "
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("E50").Value = "'  +2.84%  "
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = "'143.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.58%  "
$ws.Range("E51").Style = "Normal"
